# Generate Report for handoff
# A new handoff just occurred for file "1452fdde-cd4b-46c9-aa94-31d2e220feca" in
# both the zh-cn and de-de locales, so its "Latest Handoff Datetime" (column D,
# row 4) is refreshed to reflect the new handoff timestamp.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-17 02:39:47"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-17 02:39:56"
